# "Changed logic to checkout on platform"
#
# - The "Specific Product" checkout sheet is no longer part of the
#   automation flow, so it is removed entirely.
# - The Network sheet's login (A2) is rotated to a new network account.
# - Selection/active-cell state on a couple of the credential sheets is
#   nudged from A3 -> A2 to reflect where the automation now leaves the
#   cursor after reading the credentials.
# - "List of Products" (the last remaining sheet) becomes the active tab.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Rotate the Network sheet's username to the new platform login.
$wsNetwork = $wb.Worksheets.Item("Network")
$wsNetwork.Range("A2").Value = "newnetwork1@mailinator.com"
$wsNetwork.Range("A2").Select()

# Vendor 1 sheet cursor moves to the username cell too.
$wsVendor1 = $wb.Worksheets.Item("Vendor 1")
$wsVendor1.Range("A2").Select()

# The "Specific Product" checkout sheet is dropped from the workbook.
$wsSpecific = $wb.Worksheets.Item("Specific Product")
$wsSpecific.Delete()

# Leave "List of Products" (now the last sheet) as the active tab.
$wb.Worksheets.Item("List of Products").Activate()
